# Add team record (Wins / Losses / Ties) columns to the roster sheet.
# New columns AD:AF, header row styled like the other header cells (AC1),
# data rows 2-48 filled with the team's W/L/T totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column titles -------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold font, border, centered) from the last
# existing header cell (AC1) onto the three new header cells so they reuse
# the same cell style instead of minting a new one.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Data rows 2-48: team record for every player row ---------------------
$wins = 79
$losses = 83
$ties = 0

$lastRow = 48
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # AD
    $ws.Cells.Item($r, 31).Value = $losses  # AE
    $ws.Cells.Item($r, 32).Value = $ties    # AF
}
